# Monster commint. Neglected commiting for some time
#
# - Select the whole used range (A1:K12) on "Munka1" (was a single-row
#   selection before).
# - Add a new worksheet "Munka2" right after "Munka1" that mirrors the
#   "Munka1" layout (same header/labels, same merged cells) but with a new
#   set of duty-cycle observations in column C (rows 4-6), and leave the
#   cursor on C7.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Munka1")

# Munka1: the active selection becomes the whole used range.
$ws1.Range("A1:K12").Select()

# --- Add "Munka2" right after "Munka1" -------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Munka2"

# Header row / label column, same as Munka1.
$ws2.Range("A1").Value = " "
$ws2.Range("C1").Value = "PWM freq (Hz)"
$ws2.Range("A3").Value = "PWM duty cycle"

# Frequency header values (row 2).
$ws2.Range("C2").Value = 50
$ws2.Range("D2").Value = 100
$ws2.Range("E2").Value = 200
$ws2.Range("F2").Value = 500
$ws2.Range("G2").Value = 1000
$ws2.Range("H2").Value = 5000
$ws2.Range("I2").Value = 10000
$ws2.Range("J2").Value = 20000
$ws2.Range("K2").Value = 40000

# Duty-cycle values (column B).
$ws2.Range("B3").Value = 10
$ws2.Range("B4").Value = 20
$ws2.Range("B5").Value = 30
$ws2.Range("B6").Value = 40
$ws2.Range("B7").Value = 50
$ws2.Range("B8").Value = 60
$ws2.Range("B9").Value = 70
$ws2.Range("B10").Value = 80
$ws2.Range("B11").Value = 90
$ws2.Range("B12").Value = 100

# New observations for this sheet.
$ws2.Range("C4").Value = "akad"
$ws2.Range("C5").Value = "lassan, lánctalp néha megakad"
$ws2.Range("C6").Value = "lassan"

# Merge the label blocks exactly like on Munka1.
$ws2.Range("A1:B2").Merge()
$ws2.Range("C1:K1").Merge()
$ws2.Range("A3:A12").Merge()

# Formatting: centered (+ vertically centered for the two merged corner
# blocks), matching Munka1.
$ws2.Range("A1:B2").HorizontalAlignment = -4108
$ws2.Range("A1:B2").VerticalAlignment = -4108
$ws2.Range("C1:K1").HorizontalAlignment = -4108
$ws2.Range("A3:A12").HorizontalAlignment = -4108

# Leave the cursor on C7, as in the saved file.
$ws2.Range("C7").Select()
